$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1908315565031983
$ws.Range("C2").Value = 0.5746268656716418
$ws.Range("J2").Value = 0.0138592750533049
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.07782515991471216

# Row 3
$ws.Range("B3").Value = 0.007207207207207207
$ws.Range("C3").Value = 0.03423423423423423
$ws.Range("J3").Value = 0.01621621621621622
$ws.Range("P3").Value = 0.7603603603603604
$ws.Range("S3").Value = 0.181981981981982

# Row 6
$ws.Range("B6").Value = 0.06774668630338733
$ws.Range("D6").Value = 0.0117820324005891
$ws.Range("E6").Value = 0.001472754050073638
$ws.Range("F6").Value = 0.07658321060382917
$ws.Range("J6").Value = 0.251840942562592
$ws.Range("O6").Value = 0.01914580265095729
$ws.Range("Q6").Value = 0.1222385861561119
$ws.Range("R6").Value = 0.08100147275405008
$ws.Range("S6").Value = 0.3681885125184094

# Row 7
$ws.Range("B7").Value = 0.1077389984825493
$ws.Range("D7").Value = 0.01669195751138088
$ws.Range("F7").Value = 0.0637329286798179
$ws.Range("J7").Value = 0.1320182094081942
$ws.Range("O7").Value = 0.02731411229135053
$ws.Range("Q7").Value = 0.1714719271623672
$ws.Range("R7").Value = 0.08345978755690441
$ws.Range("S7").Value = 0.3975720789074355

# Row 8
$ws.Range("B8").Value = 0.09783368273934312
$ws.Range("D8").Value = 0.02026554856743536
$ws.Range("F8").Value = 0.05590496156533892
$ws.Range("J8").Value = 0.1118099231306778
$ws.Range("O8").Value = 0.02725366876310273
$ws.Range("Q8").Value = 0.1586303284416492
$ws.Range("R8").Value = 0.08735150244584207
$ws.Range("S8").Value = 0.4409503843466108

# Row 9
$ws.Range("B9").Value = 0.09650582362728785
$ws.Range("D9").Value = 0.009983361064891847
$ws.Range("E9").Value = 0.001663893510815308
$ws.Range("F9").Value = 0.0632279534109817
$ws.Range("J9").Value = 0.1381031613976705
$ws.Range("O9").Value = 0.01996672212978369
$ws.Range("Q9").Value = 0.1331114808652246
$ws.Range("R9").Value = 0.09317803660565724
$ws.Range("S9").Value = 0.4442595673876872

# Row 10
$ws.Range("B10").Value = 0.1079100328531716
$ws.Range("D10").Value = 0.01769016932019207
$ws.Range("E10").Value = 0.0005054334091483447
$ws.Range("F10").Value = 0.06823351023502654
$ws.Range("J10").Value = 0.1190295678544352
$ws.Range("O10").Value = 0.01541571897902451
$ws.Range("Q10").Value = 0.2158200657063432
$ws.Range("R10").Value = 0.08011119535001264
$ws.Range("S10").Value = 0.3752843062926459

# Row 11
$ws.Range("G11").Value = 0.1472134595162986
$ws.Range("J11").Value = 0.06624605678233439
$ws.Range("K11").Value = 0.195583596214511
$ws.Range("L11").Value = 0.573080967402734
$ws.Range("S11").Value = 0.01787592008412198

# Row 12
$ws.Range("G12").Value = 0.7544169611307421
$ws.Range("J12").Value = 0.176678445229682
$ws.Range("K12").Value = 0.00530035335689046
$ws.Range("L12").Value = 0.03003533568904593
$ws.Range("S12").Value = 0.03356890459363958

# Row 14
$ws.Range("G14").Value = 0.625
$ws.Range("J14").Value = 0.25
$ws.Range("S14").Value = 0.125

# Row 15
$ws.Range("F15").Value = 0.01047904191616766
$ws.Range("H15").Value = 0.1616766467065868
$ws.Range("I15").Value = 0.05389221556886228
$ws.Range("J15").Value = 0.3592814371257485
$ws.Range("K15").Value = 0.05838323353293413
$ws.Range("M15").Value = 0.01197604790419162
$ws.Range("N15").Value = 0.002994011976047904
$ws.Range("O15").Value = 0.0658682634730539
$ws.Range("S15").Value = 0.2754491017964072

# Row 16
$ws.Range("F16").Value = 0.0131578947368421
$ws.Range("H16").Value = 0.1825657894736842
$ws.Range("I16").Value = 0.08388157894736842
$ws.Range("J16").Value = 0.4029605263157895
$ws.Range("K16").Value = 0.1282894736842105
$ws.Range("M16").Value = 0.01644736842105263
$ws.Range("N16").Value = 0.001644736842105263
$ws.Range("O16").Value = 0.04111842105263158
$ws.Range("S16").Value = 0.1299342105263158

# Row 17
$ws.Range("F17").Value = 0.01783060921248143
$ws.Range("H17").Value = 0.175334323922734
$ws.Range("I17").Value = 0.08989598811292719
$ws.Range("J17").Value = 0.4331352154531947
$ws.Range("K17").Value = 0.1040118870728083
$ws.Range("M17").Value = 0.02451708766716196
$ws.Range("N17").Value = 0.002228826151560178
$ws.Range("O17").Value = 0.05943536404160475
$ws.Range("S17").Value = 0.09361069836552749

# Row 18
$ws.Range("F18").Value = 0.02003338898163606
$ws.Range("H18").Value = 0.1886477462437396
$ws.Range("I18").Value = 0.08848080133555926
$ws.Range("J18").Value = 0.4190317195325542
$ws.Range("K18").Value = 0.1085141903171953
$ws.Range("M18").Value = 0.01335559265442404
$ws.Range("O18").Value = 0.05843071786310518
$ws.Range("S18").Value = 0.1035058430717863

# Row 19
$ws.Range("F19").Value = 0.0171414166876733
$ws.Range("H19").Value = 0.217544744139148
$ws.Range("I19").Value = 0.08721956138139653
$ws.Range("J19").Value = 0.375346609528611
$ws.Range("K19").Value = 0.1106629694983615
$ws.Range("M19").Value = 0.02268716914544996
$ws.Range("N19").Value = 0.001008318628686665
$ws.Range("O19").Value = 0.0642803125787749
$ws.Range("S19").Value = 0.1041088984118982
